# Plantilla de Casos de Uso - actualización de contenido
# Adds the 7 new Caso de Uso descriptions (cuenta, sesión, contraseña,
# personal de evento) and re-points rows 42-48 ("Crear cuenta" ...
# "Eliminar personal del evento") at their specific description text
# instead of the generic placeholder. Also bumps the priority of the
# last two new use cases (Agregar/Eliminar personal del evento) from
# "Baja" to "Alta", and updates the sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- New, specific descriptions for column C (replacing the generic
# placeholder text previously shared by rows 42-48) -----------------
$ws.Range("C42").Value = "El Usuario deberá poder crear una cuenta en el sistema para poder acceder a las funciones de la misma de manera sencilla e intuitiva.`n"
$ws.Range("C43").Value = "El usuario deberá poder modificar la información de su cuenta para que tenga control sobre los datos que están ahí de manera sencilla e intuitiva.`n"
$ws.Range("C44").Value = "El usuario deberá poder iniciar sesión en el sistema para que exista seguridad y privacidad en sus datos de manera sencilla, intuitiva y rápida.`n"
$ws.Range("C45").Value = "El usuario  deberá poder cerrar sesión en el sistema para tener un control de la seguridad y privacidad de su información de manera sencilla e intuitiva.`n"
$ws.Range("C46").Value = "El usuario deberá poder recuperar su contraseña para tener control de su cuenta y poder recuperarla en caso de olvidar su contraseña esto de manera segura e intuitiva.`n"
$ws.Range("C47").Value = "El líder del evento deberá poder agregar personal al evento para que puedan participar en las activiades del mismo de manera sencilla e intuitiva.`n"
$ws.Range("C48").Value = "El líder del evento deberá poder eliminar personal del evento para poder tener control de quienes ya no participan en el evento  esto de manera sencilla e intuitiva.`n"

# --- Priority bump: "Agregar personal al evento" / "Eliminar personal
# del evento" move from Baja to Alta -------------------------------
$ws.Range("H47").Value = "Alta"
$ws.Range("H48").Value = "Alta"

# --- Sheet view / selection state -----------------------------------
$ws.Activate()
$ws.Range("D29:D48").Select()
$excel.ActiveWindow.Zoom = 91
